$wb = $excel.ActiveWorkbook

# Reference sheet holding the existing "Menlo" direct formatting we want to
# reuse for the new brand list (avoids creating duplicate font/style entries).
$sourceSheet = $wb.Worksheets.Item(1)

# Add the new "TradeInDevice" worksheet as the last sheet in the workbook.
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "TradeInDevice"

# Populate the trade-in brand list (header + brand names).
$values = @("brands", "Apple", "Samsung", "Google", "LG", "Motorola", "HTC", "OnePlus", "Microsoft", "Nokia", "Other")
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 2 + $i
    $newSheet.Cells.Item($row, 1).Value = $values[$i]
}

# Match the monospace formatting used elsewhere in the workbook (copy the
# direct formatting from an existing "Menlo" cell onto the brand rows).
$sourceSheet.Range("A2").Copy()
$newSheet.Range("A3:A12").PasteSpecial(-4122)

$newSheet.Columns.Item(1).ColumnWidth = 13.33

# Select D14 on the new sheet and make it the active tab (matches the saved
# view state).
$newSheet.Range("D14").Select()
